$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New time-sheet entry for row 30: date, hours, and activity description.
# Column C (running total) is a formula already present in the sheet and will
# recalculate automatically once B30 receives a value.
$ws.Range("A30").Value = 44248
$ws.Range("A30").NumberFormat = "mm/dd/yy"
$ws.Range("B30").Value = 8
$ws.Range("D30").Value = "Implementierung UI-Logik für die Decoder-Channel-Auswahl"

# Update the view so it matches where the user left off after the edit.
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D24").Select()
